$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Förändrad" (changed) date column C bumped by one day (45180 -> 45181)
# for every data row (row 2 through row 319).
for ($r = 2; $r -le 319; $r++) {
    $ws.Cells.Item($r, 3).Value = 45181
}

# Row 3 ("A 33191-2022") gained an extra signal species: "Gul taggsvamp".
# That bumps NT (col J), Rödlistade (col O) and Alla arter (col Q) counts,
# and the species list in column R gets the new line inserted after "Knärot".
$ws.Cells.Item(3, 10).Value = 2
$ws.Cells.Item(3, 15).Value = 3
$ws.Cells.Item(3, 17).Value = 5
$ws.Cells.Item(3, 18).Value = "Knärot`r`nGul taggsvamp`r`nSkirmossa`r`nHavstulpanlav`r`nKorallblylav"
